# Update "想去人数" (want-to-go count) figures in column F
# on the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 60
$ws1.Range("F4").Value  = 79
$ws1.Range("F7").Value  = 2645
$ws1.Range("F8").Value  = 1156
$ws1.Range("F9").Value  = 237
$ws1.Range("F10").Value = 96
$ws1.Range("F11").Value = 7634
$ws1.Range("F13").Value = 238
$ws1.Range("F14").Value = 588
$ws1.Range("F15").Value = 11627
$ws1.Range("F16").Value = 11884
$ws1.Range("F18").Value = 81
$ws1.Range("F20").Value = 18
$ws1.Range("F21").Value = 64

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 60
$ws4.Range("F4").Value  = 79
$ws4.Range("F7").Value  = 2645
$ws4.Range("F9").Value  = 1156
$ws4.Range("F10").Value = 237
$ws4.Range("F11").Value = 96
$ws4.Range("F12").Value = 7640
$ws4.Range("F14").Value = 238
$ws4.Range("F15").Value = 588
$ws4.Range("F16").Value = 11627
$ws4.Range("F17").Value = 11884
$ws4.Range("F19").Value = 81
$ws4.Range("F21").Value = 18
$ws4.Range("F22").Value = 64
